$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text even for numeric-looking strings,
# so Excel does not auto-convert values like "481.03" into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '68.618.18'
$ws.Range('E2').Value = '  +2.15%  '
$ws.Range('D3').Value = '3.918.66'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '481.03'
$ws.Range('D6').Value = '144.36'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').Value = '0.167'
$ws.Range('E10').Value = '  +8.48%  '
$ws.Range('E11').Value = '  +12.09%  '
$ws.Range('D12').Value = '42.48'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = '10.46'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '4.548.15'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').Value = '14.57'
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').Value = '3.899.18'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').Value = '19.61'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('E19').Value = '  -2.92%  '
$ws.Range('D20').Value = '68.690.77'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '431.76'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').Value = '14.57'
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('D24').Value = '87.21'
$ws.Range('E24').Value = '  -2.04%  '
$ws.Range('D25').Value = '11.56'
$ws.Range('E25').Value = '  +16.01%  '
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').Value = '37.93'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = '10.14'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = '5.83'
$ws.Range('E29').Value = '  +6.26%  '
$ws.Range('D30').Value = '703.29'
$ws.Range('E30').Value = '  -3.71%  '
$ws.Range('E31').Value = '  -3.23%  '
$ws.Range('E32').Value = '  -4.31%  '
$ws.Range('D34').Value = '0.0₃0895'
$ws.Range('E34').Value = '  +29.95%  '
$ws.Range('D35').Value = '41.22'
$ws.Range('E35').Value = '  -8.09%  '
$ws.Range('D36').Value = '59.00'
$ws.Range('E36').Value = '  +1.80%  '
$ws.Range('D37').Value = '0.151'
$ws.Range('E37').Value = '  -7.55%  '
$ws.Range('D38').Value = '5.64'
$ws.Range('E38').Value = '  +1.72%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('E40').Value = '  -2.47%  '
$ws.Range('E41').Value = '  +9.55%  '
$ws.Range('D42').Value = '2.71'
$ws.Range('E42').Value = '  +5.09%  '
$ws.Range('E43').Value = '  +2.09%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.141'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.338'
$ws.Range('E45').Value = '  -3.33%  '
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').Value = '3.40'
$ws.Range('E47').Value = '  -1.32%  '
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').Value = '147.25'
$ws.Range('E49').Value = '  +2.04%  '
$ws.Range('D50').Value = '3.14'
$ws.Range('E50').Value = '  -4.29%  '
$ws.Range('E51').Value = '  -2.29%  '

# Restore original (default/general) formatting on the price column so the
# saved cells carry no extra style index, matching the source workbook.
$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"
